$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking label as literal text (matches the
# original "shared string" storage) instead of letting Excel infer a
# Number type, then drop the temporary Text number-format again so the
# cell's style index is left untouched.
function Set-TextValue($rangeAddress, $text) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.NumberFormat = ""
}

# Header row: service names -> model years
Set-TextValue "B1" "1985"
Set-TextValue "C1" "1992"
Set-TextValue "D1" "1995"

# Row labels: ServiceA/B/C -> CarModelA/B/C
$ws.Range("A2").Value = "CarModelA"
$ws.Range("A3").Value = "CarModelB"
$ws.Range("A4").Value = "CarModelC"

# The "X" marker moves from C2 to (removed), C3 to (removed), B4 to (removed),
# D4 to C4 - i.e. the mapping matrix is re-diagonalized.
$ws.Range("C2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("C4").Value = "X"
